$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44789
$ws.Range("J2").Value = 90

$ws.Range("D3").Value = 44407
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = 25000
$ws.Range("P3").Value = 1667

$ws.Range("D4").Value = 44838
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 22000
$ws.Range("P4").Value = 1467

$ws.Range("D6").Value = 44781
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 24000
$ws.Range("L6").Value = 24000
$ws.Range("M6").Value = 24000
$ws.Range("P6").Value = 1600

$ws.Range("D7").Value = 44817
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 23000
$ws.Range("L7").Value = 23000
$ws.Range("M7").Value = 23000
$ws.Range("P7").Value = 1533

$ws.Range("D8").Value = 44365
$ws.Range("J8").Value = 80

$ws.Range("D10").Value = 44764
$ws.Range("J10").Value = 90
$ws.Range("K10").Value = 24000
$ws.Range("L10").Value = 24000
$ws.Range("M10").Value = 24000
$ws.Range("P10").Value = 1600

$ws.Range("D11").Value = 44827
$ws.Range("J11").Value = 90

$ws.Range("D12").Value = 44771
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 25000
$ws.Range("P12").Value = 1667

$ws.Range("D13").Value = 44803
$ws.Range("K13").Value = 24000
$ws.Range("L13").Value = 24000
$ws.Range("M13").Value = 24000
$ws.Range("P13").Value = 1600

$ws.Range("D14").Value = 44792
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 24000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 24000
$ws.Range("P14").Value = 1600

$ws.Range("D15").Value = 44754

$ws.Range("D16").Value = 44775
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 24000
$ws.Range("M16").Value = 24000
$ws.Range("P16").Value = 1600

$ws.Range("D17").Value = 44782
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 24000
$ws.Range("L17").Value = 24000
$ws.Range("M17").Value = 24000
$ws.Range("P17").Value = 1600

$ws.Range("D18").Value = 44750
$ws.Range("K18").Value = 25000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 25000
$ws.Range("P18").Value = 1667

$ws.Range("D20").Value = 44799
$ws.Range("K20").Value = 23000
$ws.Range("L20").Value = 23000
$ws.Range("M20").Value = 23000
$ws.Range("P20").Value = 1533

$ws.Range("D21").Value = 44831
$ws.Range("J21").Value = 90
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 25000
$ws.Range("P21").Value = 1667

$ws.Range("D23").Value = 44740

$ws.Range("D24").Value = 44400
$ws.Range("J24").Value = 80

$ws.Range("D25").Value = 44819
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 22000
$ws.Range("L25").Value = 22000
$ws.Range("M25").Value = 22000
$ws.Range("P25").Value = 1467

$ws.Range("D26").Value = 44757
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 25000
$ws.Range("P26").Value = 1667

$ws.Range("D27").Value = 44806
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 23000
$ws.Range("L27").Value = 23000
$ws.Range("M27").Value = 23000
$ws.Range("P27").Value = 1533

